# Update Model data base for methanol (Units sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# --- Row 2 (Solar_Plant_Kasso) ---
# fom_cost updated; several ramp/start/shutdown/min-op values cleared
$ws.Range("O2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("AB2").Value = 1.29
$ws.Range("AH2").ClearContents()
$ws.Range("AJ2").ClearContents()

# --- Row 3 (Electrolyzer) ---
$ws.Range("W3").Value = 0.0063
$ws.Range("Y3").Value = 1.76
$ws.Range("AD3").ClearContents()
$ws.Range("AB3").Value = 4.34

# --- Row 4 (CO2_Vaporizer) ---
$ws.Range("C4").Value = "Power_Kasso"
$ws.Range("W4").Value = 280.5

# --- Row 5 (Destilation_Tower) ---
$ws.Range("C5").Value = "Steam"
$ws.Range("W5").Value = 11.99

# --- Row 6 (Methanol_Reactor) ---
$ws.Range("J6").Value = 52
$ws.Range("L6").ClearContents()
$ws.Range("S6").Value = 0.5
$ws.Range("U6").Value = 0.5
$ws.Range("W6").Value = 4.57
$ws.Range("Y6").Value = 4.32
$ws.Range("AA6").ClearContents()
$ws.Range("AB6").Value = 4.45

# --- View state: active cell moved to L10, scroll position reset ---
$ws.Range("L10").Select()
